$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 1 de Junio de 2020 a las 00:05"

# --- Estados Unidos (row 4) ---
$ws.Range("B4").Value = 1834502
$ws.Range("C4").Value = 17682
$ws.Range("D4").Value = 541168
$ws.Range("E4").Value = 1187193
$ws.Range("G4").Value = 584
$ws.Range("H4").Value = 106141

# --- Brasil (row 5) ---
$ws.Range("B5").Value = 506708
$ws.Range("C5").Value = 8268
$ws.Range("E5").Value = 272236
$ws.Range("G5").Value = 267
$ws.Range("H5").Value = 29101

# --- Peru overtakes Turquia: row 13 becomes Peru, row 14 becomes Turquia ---
$ws.Range("A13").Value = "Peru"
$ws.Range("B13").Value = 164476
$ws.Range("C13").Value = 8805
$ws.Range("D13").Value = 67208
$ws.Range("E13").Value = 92762
$ws.Range("G13").Value = 135
$ws.Range("H13").Value = 4506

$ws.Range("A14").Value = "Turquia"
$ws.Range("B14").Value = 163942
$ws.Range("C14").Value = 839
$ws.Range("D14").Value = 127973
$ws.Range("E14").Value = 31429
$ws.Range("G14").Value = 25
$ws.Range("H14").Value = 4540

# --- Suiza (row 33) ---
$ws.Range("D33").Value = 28500
$ws.Range("E33").Value = 442

# --- Egipto (row 38) ---
$ws.Range("D38").Value = 6037
$ws.Range("E38").Value = 17989

# --- Estado de Palestina (row 142) ---
$ws.Range("D142").Value = 372
$ws.Range("E142").Value = 73

# --- Ruanda (row 147) ---
$ws.Range("B147").Value = 370
$ws.Range("C147").Value = 11
$ws.Range("D147").Value = 256
$ws.Range("E147").Value = 113

# --- Montserrat overtakes Seychelles (tie-break): row 210 becomes Montserrat, row 211 becomes Seychelles ---
$ws.Range("A210").Value = "Montserrat"
$ws.Range("D210").Value = 10
$ws.Range("H210").Value = 1

$ws.Range("A211").Value = "Seychelles"
$ws.Range("D211").Value = 11
$ws.Range("H211").Value = 0
